$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content changes -------------------------------------------------
# New row 8 (Codice 7): "Revisione onomastica stradale e dei numeri civici".
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Revisione onomastica stradale e dei numeri civici"

# Row 4 (Codice 3): text changes from "Motivo specie famiglia convivenza"
# to "Dati generali famiglia convivenza".
$ws.Range("B4").Value = "Dati generali famiglia convivenza"

# --- Column width -----------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 51.85546875

# --- Selection ---------------------------------------------------------------
$ws.Range("B5").Select()

# --- Borders -----------------------------------------------------------------
# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
# xlThin=2, xlMedium=-4138
# xlContinuous=1

function Set-Edge($rng, $edge, $weight) {
    $b = $rng.Borders.Item($edge)
    $b.LineStyle = 1
    $b.Weight = $weight
}

$thin = 2
$medium = -4138

# Row 5: A5 thin L/R/Top; B5 medium L/R
Set-Edge $ws.Range("A5") 7 $thin
Set-Edge $ws.Range("A5") 10 $thin
Set-Edge $ws.Range("A5") 8 $thin
Set-Edge $ws.Range("B5") 7 $medium
Set-Edge $ws.Range("B5") 10 $medium

# Row 6: A6 medium L, thin R, medium Top; B6 medium L/R/Top
Set-Edge $ws.Range("A6") 7 $medium
Set-Edge $ws.Range("A6") 10 $thin
Set-Edge $ws.Range("A6") 8 $medium
Set-Edge $ws.Range("B6") 7 $medium
Set-Edge $ws.Range("B6") 10 $medium
Set-Edge $ws.Range("B6") 8 $medium

# Row 7: A7 medium L, thin R, thin Top, thin Bottom; B7 medium box (all sides)
Set-Edge $ws.Range("A7") 7 $medium
Set-Edge $ws.Range("A7") 10 $thin
Set-Edge $ws.Range("A7") 8 $thin
Set-Edge $ws.Range("A7") 9 $thin
Set-Edge $ws.Range("B7") 7 $medium
Set-Edge $ws.Range("B7") 10 $medium
Set-Edge $ws.Range("B7") 8 $medium
Set-Edge $ws.Range("B7") 9 $medium

# Row 8: A8 medium L, thin R, thin Top, medium Bottom; B8 thin L, medium R, thin Top, medium Bottom
Set-Edge $ws.Range("A8") 7 $medium
Set-Edge $ws.Range("A8") 10 $thin
Set-Edge $ws.Range("A8") 8 $thin
Set-Edge $ws.Range("A8") 9 $medium
Set-Edge $ws.Range("B8") 7 $thin
Set-Edge $ws.Range("B8") 10 $medium
Set-Edge $ws.Range("B8") 8 $thin
Set-Edge $ws.Range("B8") 9 $medium
